$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28..148 down to 29..149
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new data record
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 45251
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100114007
$ws.Range("G28").Value = "Jengibre"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 430
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 14000
$ws.Range("M28").Value = 14000
$ws.Range("N28").Value = "`$/caja 15 kilos"
$ws.Range("O28").Value = "Perú"
$ws.Range("P28").Value = 933
$ws.Range("Q28").Value = 15
$ws.Range("R28").Value = "Hortaliza"
